$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.670.76"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "1.565.56"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.27"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "25.21"
$ws.Range("E8").Value = "  +5.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.244"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0894"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").Value = "1.788.49"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "1.573.16"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").Value = "28.676.99"
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.23"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.82"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.997"
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.05"
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.62"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.78"
$ws.Range("E26").Value = "  -0.74%  "
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0460"
$ws.Range("E30").Value = "  -3.97%  "
$ws.Range("E31").Value = "  -2.27%  "
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").Value = "1.393.42"
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("E34").Value = "  -3.14%  "
$ws.Range("E35").Value = "  -4.41%  "
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.29"
$ws.Range("E38").Value = "  -2.41%  "
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("E40").Value = "  +1.24%  "
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.771"
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("E44").Value = "  -3.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.03"
$ws.Range("E45").Value = "  +2.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.24"
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("D47").Value = "1.701.82"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.868"
$ws.Range("E48").Value = "  -5.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.11"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "43.27"
$ws.Range("E50").Value = "  +6.83%  "
$ws.Range("E51").Value = "  -0.66%  "
